$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2024-12-20 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-12-21 Saturday", 2) | Out-Null

# Update each arithmetic-answer cell in the table by (row, col) position
# so that duplicate old values at different positions map to their own new values.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "65-26=39"
$t.Cell(1, 2).Range.Text = "87+4=91"
$t.Cell(1, 3).Range.Text = "86-59=27"
$t.Cell(1, 4).Range.Text = "49+48=97"
$t.Cell(1, 5).Range.Text = "95-7=88"
$t.Cell(2, 1).Range.Text = "38+4=42"
$t.Cell(2, 2).Range.Text = "48+29=77"
$t.Cell(2, 3).Range.Text = "27-18=9"
$t.Cell(2, 4).Range.Text = "90-6=84"
$t.Cell(2, 5).Range.Text = "19+45=64"
$t.Cell(3, 1).Range.Text = "61-13=48"
$t.Cell(3, 2).Range.Text = "77-29=48"
$t.Cell(3, 3).Range.Text = "94-39=55"
$t.Cell(3, 4).Range.Text = "88+6=94"
$t.Cell(3, 5).Range.Text = "82-34=48"
$t.Cell(4, 1).Range.Text = "82-46=36"
$t.Cell(4, 2).Range.Text = "61-37=24"
$t.Cell(4, 3).Range.Text = "87-38=49"
$t.Cell(4, 4).Range.Text = "40-28=12"
$t.Cell(4, 5).Range.Text = "39+9=48"
$t.Cell(5, 1).Range.Text = "12-4=8"
$t.Cell(5, 2).Range.Text = "45-8=37"
$t.Cell(5, 3).Range.Text = "35-19=16"
$t.Cell(5, 4).Range.Text = "90-33=57"
$t.Cell(5, 5).Range.Text = "76+16=92"
$t.Cell(6, 1).Range.Text = "95-89=6"
$t.Cell(6, 2).Range.Text = "30-3=27"
$t.Cell(6, 3).Range.Text = "29+52=81"
$t.Cell(6, 4).Range.Text = "47-8=39"
$t.Cell(6, 5).Range.Text = "35+9=44"
$t.Cell(7, 1).Range.Text = "37+37=74"
$t.Cell(7, 2).Range.Text = "58-19=39"
$t.Cell(7, 3).Range.Text = "17+9=26"
$t.Cell(7, 4).Range.Text = "19+14=33"
$t.Cell(7, 5).Range.Text = "48+48=96"
$t.Cell(8, 1).Range.Text = "19+26=45"
$t.Cell(8, 2).Range.Text = "9+66=75"
$t.Cell(8, 3).Range.Text = "69+19=88"
$t.Cell(8, 4).Range.Text = "28+67=95"
$t.Cell(8, 5).Range.Text = "62-55=7"
$t.Cell(9, 1).Range.Text = "31-29=2"
$t.Cell(9, 2).Range.Text = "37+18=55"
$t.Cell(9, 3).Range.Text = "66+6=72"
$t.Cell(9, 5).Range.Text = "9+25=34"
$t.Cell(10, 1).Range.Text = "59+15=74"
$t.Cell(10, 2).Range.Text = "38-29=9"
$t.Cell(10, 3).Range.Text = "94-57=37"
$t.Cell(10, 4).Range.Text = "34-18=16"
$t.Cell(10, 5).Range.Text = "28+59=87"
$t.Cell(11, 1).Range.Text = "91-53=38"
$t.Cell(11, 2).Range.Text = "40-1=39"
$t.Cell(11, 3).Range.Text = "46-7=39"
$t.Cell(11, 4).Range.Text = "65-26=39"
$t.Cell(11, 5).Range.Text = "25+17=42"
$t.Cell(12, 1).Range.Text = "47+48=95"
$t.Cell(12, 2).Range.Text = "44+39=83"
$t.Cell(12, 3).Range.Text = "5+57=62"
$t.Cell(12, 4).Range.Text = "9+39=48"
$t.Cell(12, 5).Range.Text = "7+77=84"
$t.Cell(13, 1).Range.Text = "16-9=7"
$t.Cell(13, 2).Range.Text = "9+5=14"
$t.Cell(13, 3).Range.Text = "24+28=52"
$t.Cell(13, 4).Range.Text = "53-14=39"
$t.Cell(13, 5).Range.Text = "18+54=72"
$t.Cell(14, 1).Range.Text = "5+89=94"
$t.Cell(14, 2).Range.Text = "7+78=85"
$t.Cell(14, 3).Range.Text = "6+69=75"
$t.Cell(14, 4).Range.Text = "19+63=82"
$t.Cell(14, 5).Range.Text = "70-26=44"
$t.Cell(15, 1).Range.Text = "83-14=69"
$t.Cell(15, 2).Range.Text = "83-55=28"
$t.Cell(15, 3).Range.Text = "79+18=97"
$t.Cell(15, 4).Range.Text = "17+64=81"
$t.Cell(15, 5).Range.Text = "92-68=24"
$t.Cell(16, 1).Range.Text = "16+79=95"
$t.Cell(16, 2).Range.Text = "83-28=55"
$t.Cell(16, 3).Range.Text = "94-37=57"
$t.Cell(16, 4).Range.Text = "27+19=46"
$t.Cell(16, 5).Range.Text = "73-29=44"
$t.Cell(17, 1).Range.Text = "76-47=29"
$t.Cell(17, 2).Range.Text = "19+8=27"
$t.Cell(17, 3).Range.Text = "77+6=83"
$t.Cell(17, 4).Range.Text = "8+28=36"
$t.Cell(17, 5).Range.Text = "65-38=27"
$t.Cell(18, 1).Range.Text = "49+12=61"
$t.Cell(18, 2).Range.Text = "62-7=55"
$t.Cell(18, 3).Range.Text = "16+29=45"
$t.Cell(18, 4).Range.Text = "91-7=84"
$t.Cell(18, 5).Range.Text = "19+35=54"
$t.Cell(19, 1).Range.Text = "79+13=92"
$t.Cell(19, 2).Range.Text = "89+5=94"
$t.Cell(19, 3).Range.Text = "74-45=29"
$t.Cell(19, 4).Range.Text = "38+13=51"
$t.Cell(19, 5).Range.Text = "81-77=4"
$t.Cell(20, 1).Range.Text = "91-82=9"
$t.Cell(20, 2).Range.Text = "9+83=92"
$t.Cell(20, 3).Range.Text = "92-14=78"
$t.Cell(20, 4).Range.Text = "24+49=73"
$t.Cell(20, 5).Range.Text = "74-17=57"
